$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').Value = '27.107.86'
$ws.Range('E2').Value = '  +1.24%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').Value = '1.569.24'
$ws.Range('E3').Value = '  +2.06%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('E4').Value = '  +0.04%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.84'
$ws.Range('E5').Value = '  +1.28%  '

# Row 6: 'XRP' -> 'XRP'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.491'
$ws.Range('E6').Value = '  +1.23%  '

# Row 7: 'USDC' -> 'USDC'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.20%  '

# Row 8: 'Solana' -> 'Solana'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.15'
$ws.Range('E8').Value = '  +4.18%  '

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.249'
$ws.Range('E9').Value = '  +1.59%  '

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0589'
$ws.Range('E10').Value = '  +1.42%  '

# Row 11: 'TRON' -> 'TRON'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0858'
$ws.Range('E11').Value = '  +0.45%  '

# Row 12: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range('D12').Value = '1.778.80'
$ws.Range('E12').Value = '  +1.31%  '

# Row 13: 'WrappedEther' -> 'WrappedEther'
$ws.Range('D13').Value = '1.582.45'
$ws.Range('E13').Value = '  +2.80%  '

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.77'
$ws.Range('E14').Value = '  +2.84%  '

# Row 15: 'Polygon' -> 'Polygon'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  +3.25%  '

# Row 16: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range('D16').Value = '27.080.56'
$ws.Range('E16').Value = '  +1.16%  '

# Row 17: 'Litecoin' -> 'Litecoin'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.03'
$ws.Range('E17').Value = '  +1.81%  '

# Row 18: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '219.21'
$ws.Range('E18').Value = '  +2.88%  '

# Row 19: 'ShibaInu' -> 'ShibaInu'
$ws.Range('D19').Value = '0.0₃0695'
$ws.Range('E19').Value = '  +2.11%  '

# Row 20: 'Chainlink' -> 'Chainlink'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.32'
$ws.Range('E20').Value = '  +1.30%  '

# Row 21: 'Dai' -> 'Dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +0.28%  '

# Row 22: 'Uniswap' -> 'Uniswap'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.09'
$ws.Range('E22').Value = '  +1.85%  '

# Row 23: 'Avalanche' -> 'Avalanche'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('E23').Value = '  +1.92%  '

# Row 24: 'Toncoin' -> 'Toncoin'
$ws.Range('E24').Value = '  +1.30%  '

# Row 25: 'Monero' -> 'Monero'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.32'
$ws.Range('E25').Value = '  +1.85%  '

# Row 26: 'Cosmos' -> 'Cosmos'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.64'
$ws.Range('E26').Value = '  +1.09%  '

# Row 27: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.01'
$ws.Range('E27').Value = '  +1.44%  '

# Row 28: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.19%  '

# Row 29: 'Stellar' -> 'Stellar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.104'
$ws.Range('E29').Value = '  +1.67%  '

# Row 30: 'Hedera' -> 'Hedera'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0472'
$ws.Range('E30').Value = '  +3.39%  '

# Row 31: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range('E31').Value = '  +0.61%  '

# Row 32: 'Filecoin' -> 'Filecoin'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.24'
$ws.Range('E32').Value = '  +0.60%  '

# Row 33: 'Maker' -> 'Maker'
$ws.Range('D33').Value = '1.460.09'
$ws.Range('E33').Value = '  +7.06%  '

# Row 34: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.07'
$ws.Range('E34').Value = '  +4.87%  '

# Row 35: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.57'
$ws.Range('E35').Value = '  +4.94%  '

# Row 36: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.970'
$ws.Range('E36').Value = '  +0.54%  '

# Row 37: 'HuobiToken' -> 'HuobiToken'
$ws.Range('E37').Value = '  +0.77%  '

# Row 38: 'VeChain' -> 'VeChain'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0165'
$ws.Range('E38').Value = '  +0.49%  '

# Row 39: 'ImmutableX' -> 'ImmutableX'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.524'
$ws.Range('E39').Value = '  +0.75%  '

# Row 40: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.817'
$ws.Range('E40').Value = '  +1.87%  '

# Row 41: 'PaxDollar' -> 'FraxShare'
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.74'
$ws.Range('E41').Value = '  +0.14%  '

# Row 42: 'FraxShare' -> 'PaxDollar'
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.14%  '

# Row 43: 'WEMIXToken' -> 'MXToken'
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.29'
$ws.Range('E43').Value = '  +4.00%  '

# Row 44: 'MXToken' -> 'WEMIXToken'
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.986'
$ws.Range('E44').Value = '  -0.93%  '

# Row 45: 'Aave' -> 'Aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.63'
$ws.Range('E45').Value = '  +2.78%  '

# Row 46: 'RenderToken' -> 'RenderToken'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.76'
$ws.Range('E46').Value = '  +2.55%  '

# Row 47: 'RocketPoolETH' -> 'RocketPoolETH'
$ws.Range('D47').Value = '1.695.92'
$ws.Range('E47').Value = '  +1.53%  '

# Row 48: 'Quant' -> 'Quant'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.44'
$ws.Range('E48').Value = '  +1.38%  '

# Row 49: 'Cronos' -> 'Cronos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0525'
$ws.Range('E49').Value = '  +3.42%  '

# Row 50: 'BabyDogeCoin' -> 'Algorand'
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0966'
$ws.Range('E50').Value = '  +2.45%  '

# Row 51: 'Algorand' -> 'USDD'
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.21%  '
